# edit.ps1 - apply the diff between draft-gandhi-spring-twamp-srpm-09 revisions.
#
# Summary of changes:
#  - Slide 13, "Content Placeholder 2" shape: resize/reposition (top moves
#    down, box shrinks), and two "for" -> "of" wording tweaks.
#  - Slide 13, "Rectangle 8" shape: grows slightly taller, and the
#    "Next Hop IPv6 Address" line becomes "Destination IPv6 Address".
#  - Slide 4, "Content Placeholder 2" shape: the three runs
#    "draft-" + "gandhi" + "--spring-" collapse into a single run
#    "draft-gandhi-spring-" (also fixing the stray double dash).

# PowerPoint stores shape position/size (Left/Top/Width/Height) as single
# precision (float32) points. Helper: given a target EMU value, find a
# points value that round-trips through float32 -> EMU to that exact EMU,
# so the saved OOXML offsets/extents match exactly.
function EmuToPt {
    param([double]$Emu)
    $base = $Emu / 12700.0
    for ($i = -4000; $i -le 4000; $i++) {
        $cand = $base + ($i * 0.000001)
        $f32 = [single]$cand
        $emuCalc = [math]::Floor([double]$f32 * 12700.0)
        if ($emuCalc -eq $Emu) {
            return $cand
        }
    }
    return $base
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 13
# ---------------------------------------------------------------------
$s13 = $p.Slides.Item(13)

# "Content Placeholder 2" shape (off/ext + two word tweaks)
$probeBox = $s13.Shapes.Item(4)

$probeBox.Top = EmuToPt 1311623
$probeBox.Height = EmuToPt 2022127

$tr = $probeBox.TextFrame2.TextRange
$full = $tr.Text

$needle1 = "MPLS label stack for SR-MPLS Policy"
$idx1 = $full.IndexOf($needle1) + 1
$tr.Characters($idx1, $needle1.Length).Text = "MPLS label stack of SR-MPLS Policy"

$full = $tr.Text
$needle2 = "] with Segment List for SRv6 Policy"
$idx2 = $full.IndexOf($needle2) + 1
$tr.Characters($idx2, $needle2.Length).Text = "] with Segment List of SRv6 Policy"

# "Rectangle 8" shape (height + address line wording)
$rect8 = $s13.Shapes.Item(6)
$rect8.Height = EmuToPt 3093154

$tr2 = $rect8.TextFrame2.TextRange
$full2 = $tr2.Text
$needle3 = ".  Destination IP Address = Next Hop IPv6 Address               ."
$idx3 = $full2.IndexOf($needle3) + 1
$tr2.Characters($idx3, $needle3.Length).Text = ".  Destination IP Address = Destination IPv6 Address            ."

# ---------------------------------------------------------------------
# Slide 4
# ---------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$history = $s4.Shapes.Item(2)

$tr4 = $history.TextFrame2.TextRange
$full4 = $tr4.Text
$needle4 = "draft-gandhi--spring-"
$idx4 = $full4.IndexOf($needle4) + 1
$tr4.Characters($idx4, $needle4.Length).Text = "draft-gandhi-spring-"
